$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# phi0 (row 8): change from 0*PI()/180 to -2*PI()/180
$ws.Range("B8").Formula = "=-2*PI()/180"

# theta0 (row 9): change from 8*PI()/180 to 0*PI()/180
$ws.Range("B9").Formula = "=0*PI()/180"
